# On person import, catch languages we can't make sense of
# Update Vince White's (row 3) primary_language to "Englishes" and
# record the unparsed other_languages value "fr;foo,ru,spa" in column N.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# J3 = primary_language: English -> Englishes
$ws.Range("J3").Value = "Englishes"

# N3 = other_languages (new cell/value)
$ws.Range("N3").Value = "fr;foo,ru,spa"

# Update the active selection to reflect where the user ended up
$ws.Range("N3").Select()
